$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "20.572.34"
$ws.Range("E2").Value = "  +1.67%  "
$ws.Range("D3").Value = "1.469.44"
$ws.Range("E4").Value = "  +0.28%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.9599"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.41%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "276.58"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.59%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3548"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.37%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3055"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.75%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.079"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +5.94%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.26"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.30%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06618"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.03%  "
$ws.Range("E12").Value = "  +0.48%  "
$ws.Range("E13").Value = "  +2.49%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "18.04"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.53%  "
$ws.Range("E15").Value = "  +2.38%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9601"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.65%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001016"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.89%  "
$ws.Range("D18").Value = "1.470.32"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.05954"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +6.15%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "68.83"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.87%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.469"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.62%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "14.42"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.82%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.17"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.25%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.273"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.58%  "
$ws.Range("D25").Value = "20.577.18"
$ws.Range("E25").Value = "  +1.52%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "145.37"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.47%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.082"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.85%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.07"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.11%  "
$ws.Range("D29").Value = "1.630.39"
$ws.Range("E29").Value = "  +2.25%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "114.32"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.66%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.953"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.37%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.905"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.80%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.07940"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.69%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7914"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.38%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.224"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +7.97%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.441"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.05%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.05683"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.27%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.699"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.42%  "
$ws.Range("E39").Value = "  +3.62%  "
$ws.Range("E40").Value = "  +2.23%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "10.22"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.82%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1845"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.37%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.262"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.75%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.512"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.01%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5204"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.47%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.03"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.93%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "119.98"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.11%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5147"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.11%  "
$ws.Range("E49").Value = "  +4.07%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06429"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.11%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.9935"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.22%  "
